$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Value edits -----------------------------------------------------
# B2: 1 -> 11
$ws.Range("B2").Value = 11
# B4: 3 -> 1
$ws.Range("B4").Value = 1

# --- Fill colours (Interior.Color sets the fg colour, Interior.PatternColor
#     sets the bg colour - both are written so the solid fill carries an
#     explicit background colour matching the foreground, like the target
#     file) -------------------------------------------------------------
# A1: solid red fill EE1111
$ws.Range("A1").Interior.Color = 0x1111EE
$ws.Range("A1").Interior.PatternColor = 0x1111EE

# B2: solid pale green fill 98FB98
$ws.Range("B2").Interior.Color = 0x98FB98
$ws.Range("B2").Interior.PatternColor = 0x98FB98

# B4: solid tomato fill FF6347
$ws.Range("B4").Interior.Color = 0x4763FF
$ws.Range("B4").Interior.PatternColor = 0x4763FF

# --- Selection ---------------------------------------------------------
$ws.Range("H14").Select() | Out-Null
